# Auto-generated Excel COM-interop script
# Applies numeric cell updates (and a few cell clears/additions) to match the target diff
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR (Leve profit-tracking data)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value2 = 90.25
$ws.Range("I11").Value2 = 90.25
$ws.Range("K11").Value2 = 90.25
$ws.Range("M11").Value2 = 49.75
$ws.Range("H19").Value2 = 1426.3158
$ws.Range("I19").Value2 = 742
$ws.Range("J19").Value2 = 1924
$ws.Range("K19").Value2 = 742
$ws.Range("L19").Value2 = 1924
$ws.Range("M19").Value2 = -567
$ws.Range("N19").Value2 = -2274
$ws.Range("H34").Value2 = 42799.8
$ws.Range("I34").Value2 = 25000
$ws.Range("J34").Value2 = 47249.75
$ws.Range("K34").Value2 = 25000
$ws.Range("L34").Value2 = 47249.75
$ws.Range("M34").Value2 = -24797
$ws.Range("N34").Value2 = -47655.75
$ws.Range("H36").Value2 = 42799.8
$ws.Range("I36").Value2 = 25000
$ws.Range("J36").Value2 = 47249.75
$ws.Range("K36").Value2 = 25000
$ws.Range("L36").Value2 = 47249.75
$ws.Range("M36").Value2 = -24285
$ws.Range("N36").Value2 = -48679.75
$ws.Range("H40").Value2 = 5149782
$ws.Range("I40").Value2 = 21950.334
$ws.Range("K40").Value2 = 21950.334
$ws.Range("M40").Value2 = -21775.334
$ws.Range("H41").Value2 = 1269.762
$ws.Range("I41").Value2 = 486.14285
$ws.Range("J41").Value2 = 1661.5714
$ws.Range("K41").Value2 = 486.14285
$ws.Range("L41").Value2 = 1661.5714
$ws.Range("M41").Value2 = -46.14285000000001
$ws.Range("N41").Value2 = -2541.5714
$ws.Range("H42").Value2 = 418.58334
$ws.Range("I42").Value2 = 68
$ws.Range("J42").Value2 = 669
$ws.Range("K42").Value2 = 204
$ws.Range("L42").Value2 = 2007
$ws.Range("M42").Value2 = 26
$ws.Range("N42").Value2 = -2467
$ws.Range("H55").Value2 = 579.6
$ws.Range("J55").Value2 = 1001
$ws.Range("L55").Value2 = 1001
$ws.Range("N55").Value2 = -1429
$ws.Range("H62").Value2 = 4569.5713
$ws.Range("I62").Value2 = 3595
$ws.Range("J62").Value2 = 4959.4
$ws.Range("K62").Value2 = 3595
$ws.Range("L62").Value2 = 4959.4
$ws.Range("M62").Value2 = -2971
$ws.Range("N62").Value2 = -6207.4
$ws.Range("H65").Value2 = 4569.5713
$ws.Range("I65").Value2 = 3595
$ws.Range("J65").Value2 = 4959.4
$ws.Range("K65").Value2 = 17975
$ws.Range("L65").Value2 = 24797
$ws.Range("M65").Value2 = -14855
$ws.Range("N65").Value2 = -31037
$ws.Range("H70").Value2 = 2000
$ws.Range("J70").Value2 = 2000
$ws.Range("L70").Value2 = 6000
$ws.Range("N70").Value2 = -6540
$ws.Range("H73").Value2 = 2000
$ws.Range("J73").Value2 = 2000
$ws.Range("L73").Value2 = 6000
$ws.Range("N73").Value2 = -7872
$ws.Range("H74").Value2 = 125005000
$ws.Range("I74").Value2 = 250005000
$ws.Range("K74").Value2 = 250005000
$ws.Range("M74").Value2 = -250004064
$ws.Range("H76").Value2 = 3979.8
$ws.Range("I76").Value2 = 4300
$ws.Range("J76").Value2 = 3499.5
$ws.Range("K76").Value2 = 4300
$ws.Range("L76").Value2 = 3499.5
$ws.Range("M76").Value2 = -3985
$ws.Range("N76").Value2 = -4129.5
$ws.Range("H77").Value2 = 125005000
$ws.Range("I77").Value2 = 250005000
$ws.Range("K77").Value2 = 1250025000
$ws.Range("M77").Value2 = -1250020320
$ws.Range("H79").Value2 = 3979.8
$ws.Range("I79").Value2 = 4300
$ws.Range("J79").Value2 = 3499.5
$ws.Range("K79").Value2 = 4300
$ws.Range("L79").Value2 = 3499.5
$ws.Range("M79").Value2 = -3208
$ws.Range("N79").Value2 = -5683.5
$ws.Range("H80").Value2 = 19231264
$ws.Range("I80").Value2 = 31250162
$ws.Range("K80").Value2 = 93750486
$ws.Range("M80").Value2 = -93749488
$ws.Range("H83").Value2 = 19231264
$ws.Range("I83").Value2 = 31250162
$ws.Range("K83").Value2 = 281251458
$ws.Range("M83").Value2 = -281246466
$ws.Range("H86").Value2 = 114702376
$ws.Range("I86").Value2 = 147470050
$ws.Range("K86").Value2 = 147470050
$ws.Range("M86").Value2 = -147468927
$ws.Range("H89").Value2 = 114702376
$ws.Range("I89").Value2 = 147470050
$ws.Range("K89").Value2 = 737350250
$ws.Range("M89").Value2 = -737344634
$ws.Range("H93").Value2 = 50000
$ws.Range("J93").Value2 = 50000
$ws.Range("L93").Value2 = 50000
$ws.Range("N93").Value2 = -54992
$ws.Range("H98").Value2 = 5125.5557
$ws.Range("I98").Value2 = 2688.3333
$ws.Range("K98").Value2 = 2688.3333
$ws.Range("M98").Value2 = -1190.3333
$ws.Range("H106").Value2 = 47626904
$ws.Range("I106").Value2 = 111116776
$ws.Range("K106").Value2 = 111116776
$ws.Range("M106").Value2 = -111116145
$ws.Range("H116").Value2 = 1130432
$ws.Range("J116").Value2 = 5388.75
$ws.Range("L116").Value2 = 5388.75
$ws.Range("N116").Value2 = -12272.75
$ws.Range("H122").Value2 = 5125.5557
$ws.Range("I122").Value2 = 2688.3333
$ws.Range("K122").Value2 = 8064.999899999999
$ws.Range("M122").Value2 = -5614.999899999999
$ws.Range("H123").Value2 = 92605.28
$ws.Range("J123").Value2 = 92605.28
$ws.Range("L123").Value2 = 92605.28
$ws.Range("N123").Value2 = -102405.28
$ws.Range("H132").Value2 = 212436.94
$ws.Range("I132").Value2 = 244806.25
$ws.Range("K132").Value2 = 734418.75
$ws.Range("M132").Value2 = -731888.75
$ws.Range("H137").Value2 = 5707.1875
$ws.Range("I137").Value2 = 4265.6665
$ws.Range("J137").Value2 = 6572.1
$ws.Range("K137").Value2 = 12796.9995
$ws.Range("L137").Value2 = 19716.3
$ws.Range("M137").Value2 = -10246.9995
$ws.Range("N137").Value2 = -24816.3
$ws.Range("H138").Value2 = 7738.537
$ws.Range("J138").Value2 = 8577.421
$ws.Range("L138").Value2 = 25732.263
$ws.Range("N138").Value2 = -36012.263
$ws.Range("H141").Value2 = 7488.5557
$ws.Range("I141").Value2 = 8732.833000000001
$ws.Range("K141").Value2 = 26198.499
$ws.Range("M141").Value2 = -21018.499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 2114.6882
$ws.Range("I32").Value2 = 1689.6714
$ws.Range("K32").Value2 = 1689.6714
$ws.Range("M32").Value2 = -1402.6714
$ws.Range("H45").Value2 = 2421.238
$ws.Range("I45").Value2 = 1263.2667
$ws.Range("J45").Value2 = 5316.1665
$ws.Range("K45").Value2 = 1263.2667
$ws.Range("L45").Value2 = 5316.1665
$ws.Range("M45").Value2 = -886.2666999999999
$ws.Range("N45").Value2 = -6070.1665
$ws.Range("H53").Value2 = 1000
$ws.Range("J53").Value2 = 0
$ws.Range("L53").Value2 = 0
$ws.Range("N53").ClearContents()
$ws.Range("H61").Value2 = 11455.235
$ws.Range("I61").Value2 = 10134.76
$ws.Range("K61").Value2 = 10134.76
$ws.Range("M61").Value2 = -9922.76
$ws.Range("H74").Value2 = 3847.2222
$ws.Range("I74").Value2 = 908.9375
$ws.Range("K74").Value2 = 908.9375
$ws.Range("M74").Value2 = -34.9375
$ws.Range("H77").Value2 = 3847.2222
$ws.Range("I77").Value2 = 908.9375
$ws.Range("K77").Value2 = 4544.6875
$ws.Range("M77").Value2 = -176.6875
$ws.Range("H80").Value2 = 89948.60000000001
$ws.Range("J80").Value2 = 89948.60000000001
$ws.Range("L80").Value2 = 89948.60000000001
$ws.Range("N80").Value2 = -91944.60000000001
$ws.Range("H83").Value2 = 89948.60000000001
$ws.Range("J83").Value2 = 89948.60000000001
$ws.Range("L83").Value2 = 269845.8
$ws.Range("N83").Value2 = -279829.8
$ws.Range("H102").Value2 = 624283.4399999999
$ws.Range("I102").Value2 = 857917.9399999999
$ws.Range("K102").Value2 = 857917.9399999999
$ws.Range("M102").Value2 = -856295.9399999999
$ws.Range("H110").Value2 = 1076931.1
$ws.Range("I110").Value2 = 1570654.9
$ws.Range("J110").Value2 = 7196.3335
$ws.Range("K110").Value2 = 1570654.9
$ws.Range("L110").Value2 = 7196.3335
$ws.Range("M110").Value2 = -1568609.9
$ws.Range("N110").Value2 = -11286.3335
$ws.Range("H131").Value2 = 93500
$ws.Range("J131").Value2 = 93500
$ws.Range("L131").Value2 = 93500
$ws.Range("N131").Value2 = -103580
$ws.Range("H132").Value2 = 10665.448
$ws.Range("I132").Value2 = 13085.68
$ws.Range("K132").Value2 = 39257.04
$ws.Range("M132").Value2 = -36727.04
$ws.Range("H136").Value2 = 11455.235
$ws.Range("I136").Value2 = 10134.76
$ws.Range("K136").Value2 = 30404.28
$ws.Range("M136").Value2 = -27854.28

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value2 = 2500
$ws.Range("J22").Value2 = 2500
$ws.Range("L22").Value2 = 2500
$ws.Range("N22").Value2 = -2846
$ws.Range("H82").Value2 = 37885
$ws.Range("J82").Value2 = 56978.668
$ws.Range("L82").Value2 = 56978.668
$ws.Range("N82").Value2 = -57744.668
$ws.Range("H85").Value2 = 37885
$ws.Range("J85").Value2 = 56978.668
$ws.Range("L85").Value2 = 56978.668
$ws.Range("N85").Value2 = -59630.668
$ws.Range("H86").Value2 = 1645.5555
$ws.Range("I86").Value2 = 1479.5333
$ws.Range("J86").Value2 = 2475.6667
$ws.Range("K86").Value2 = 1479.5333
$ws.Range("L86").Value2 = 2475.6667
$ws.Range("M86").Value2 = -356.5333000000001
$ws.Range("N86").Value2 = -4721.6667
$ws.Range("H89").Value2 = 1645.5555
$ws.Range("I89").Value2 = 1479.5333
$ws.Range("J89").Value2 = 2475.6667
$ws.Range("K89").Value2 = 7397.6665
$ws.Range("L89").Value2 = 12378.3335
$ws.Range("M89").Value2 = -1781.6665
$ws.Range("N89").Value2 = -23610.3335
$ws.Range("H94").Value2 = 444473.34
$ws.Range("I94").Value2 = 550577.75
$ws.Range("J94").Value2 = 2371.8333
$ws.Range("K94").Value2 = 550577.75
$ws.Range("L94").Value2 = 2371.8333
$ws.Range("M94").Value2 = -550126.75
$ws.Range("N94").Value2 = -3273.8333
$ws.Range("H107").Value2 = 1340.5
$ws.Range("I107").Value2 = 1211.7
$ws.Range("K107").Value2 = 1211.7
$ws.Range("M107").Value2 = 708.3
$ws.Range("H125").Value2 = 7045
$ws.Range("J125").Value2 = 0
$ws.Range("L125").Value2 = 0
$ws.Range("N125").ClearContents()
$ws.Range("H134").Value2 = 3306.4
$ws.Range("I134").Value2 = 1915.2858
$ws.Range("J134").Value2 = 7800.769
$ws.Range("K134").Value2 = 5745.857400000001
$ws.Range("L134").Value2 = 23402.307
$ws.Range("M134").Value2 = -3210.857400000001
$ws.Range("N134").Value2 = -28472.307
$ws.Range("H137").Value2 = 0
$ws.Range("I137").Value2 = 0
$ws.Range("K137").Value2 = 0
$ws.Range("M137").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value2 = 2140.8333
$ws.Range("I16").Value2 = 1786.5
$ws.Range("J16").Value2 = 2849.5
$ws.Range("K16").Value2 = 1786.5
$ws.Range("L16").Value2 = 2849.5
$ws.Range("M16").Value2 = -1499.5
$ws.Range("N16").Value2 = -3423.5
$ws.Range("H31").Value2 = 20837666
$ws.Range("I31").Value2 = 76926184
$ws.Range("J31").Value2 = 4789.1143
$ws.Range("K31").Value2 = 76926184
$ws.Range("L31").Value2 = 4789.1143
$ws.Range("M31").Value2 = -76925889
$ws.Range("N31").Value2 = -5379.1143
$ws.Range("H34").Value2 = 20837666
$ws.Range("I34").Value2 = 76926184
$ws.Range("J34").Value2 = 4789.1143
$ws.Range("K34").Value2 = 76926184
$ws.Range("L34").Value2 = 4789.1143
$ws.Range("M34").Value2 = -76925982
$ws.Range("N34").Value2 = -5193.1143
$ws.Range("H58").Value2 = 836442.2
$ws.Range("I58").Value2 = 1252351
$ws.Range("K58").Value2 = 1252351
$ws.Range("M58").Value2 = -1252148
$ws.Range("H68").Value2 = 73105.586
$ws.Range("J68").Value2 = 75206.09
$ws.Range("L68").Value2 = 75206.09
$ws.Range("N68").Value2 = -76704.09
$ws.Range("H71").Value2 = 73105.586
$ws.Range("J71").Value2 = 75206.09
$ws.Range("L71").Value2 = 225618.27
$ws.Range("N71").Value2 = -233106.27
$ws.Range("H74").Value2 = 62730.2
$ws.Range("J74").Value2 = 62730.2
$ws.Range("L74").Value2 = 62730.2
$ws.Range("N74").Value2 = -64478.2
$ws.Range("H77").Value2 = 62730.2
$ws.Range("J77").Value2 = 62730.2
$ws.Range("L77").Value2 = 188190.6
$ws.Range("N77").Value2 = -196926.6
$ws.Range("H105").Value2 = 2842126.5
$ws.Range("I105").Value2 = 7576404
$ws.Range("J105").Value2 = 1559.8
$ws.Range("K105").Value2 = 7576404
$ws.Range("L105").Value2 = 1559.8
$ws.Range("M105").Value2 = -7574657
$ws.Range("N105").Value2 = -5053.8
$ws.Range("H113").Value2 = 2140.8333
$ws.Range("I113").Value2 = 1786.5
$ws.Range("J113").Value2 = 2849.5
$ws.Range("K113").Value2 = 1786.5
$ws.Range("L113").Value2 = 2849.5
$ws.Range("M113").Value2 = 383.5
$ws.Range("N113").Value2 = -7189.5
$ws.Range("H122").Value2 = 3229.7778
$ws.Range("J122").Value2 = 4916.909
$ws.Range("L122").Value2 = 14750.727
$ws.Range("N122").Value2 = -19650.727
$ws.Range("H132").Value2 = 38471620
$ws.Range("I132").Value2 = 47621580
$ws.Range("J132").Value2 = 41798.4
$ws.Range("K132").Value2 = 142864740
$ws.Range("L132").Value2 = 125395.2
$ws.Range("M132").Value2 = -142862210
$ws.Range("N132").Value2 = -130455.2
$ws.Range("H134").Value2 = 1599.6666
$ws.Range("I134").Value2 = 899.5
$ws.Range("J134").Value2 = 3000
$ws.Range("K134").Value2 = 2698.5
$ws.Range("L134").Value2 = 9000
$ws.Range("M134").Value2 = -163.5
$ws.Range("N134").Value2 = -14070
$ws.Range("H136").Value2 = 836442.2
$ws.Range("I136").Value2 = 1252351
$ws.Range("K136").Value2 = 3757053
$ws.Range("M136").Value2 = -3754503

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value2 = 93.28570999999999
$ws.Range("I8").Value2 = 93.28570999999999
$ws.Range("K8").Value2 = 279.85713
$ws.Range("M8").Value2 = -140.85713
$ws.Range("H12").Value2 = 179.28572
$ws.Range("J12").Value2 = 175.8
$ws.Range("L12").Value2 = 527.4000000000001
$ws.Range("N12").Value2 = -873.4000000000001
$ws.Range("H18").Value2 = 34.75
$ws.Range("I18").Value2 = 34.75
$ws.Range("K18").Value2 = 104.25
$ws.Range("M18").Value2 = 64.75
$ws.Range("H21").Value2 = 3266.2666
$ws.Range("I21").Value2 = 2995
$ws.Range("K21").Value2 = 8985
$ws.Range("M21").Value2 = -8812
$ws.Range("H61").Value2 = 1250
$ws.Range("I61").Value2 = 0
$ws.Range("J61").Value2 = 1250
$ws.Range("K61").Value2 = 0
$ws.Range("L61").Value2 = 3750
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value2 = -4180
$ws.Range("H68").Value2 = 244973.72
$ws.Range("I68").Value2 = 1799.8
$ws.Range("K68").Value2 = 5399.4
$ws.Range("M68").Value2 = -4588.4
$ws.Range("H71").Value2 = 244973.72
$ws.Range("I71").Value2 = 1799.8
$ws.Range("K71").Value2 = 16198.2
$ws.Range("M71").Value2 = -12142.2
$ws.Range("H107").Value2 = 1060493.2
$ws.Range("I107").Value2 = 2293.5
$ws.Range("J107").Value2 = 1765959.6
$ws.Range("K107").Value2 = 6880.5
$ws.Range("L107").Value2 = 5297878.800000001
$ws.Range("M107").Value2 = -4960.5
$ws.Range("N107").Value2 = -5301718.800000001
$ws.Range("H122").Value2 = 849.2857
$ws.Range("I122").Value2 = 690.3077
$ws.Range("K122").Value2 = 6212.7693
$ws.Range("M122").Value2 = -3762.7693

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value2 = 21154400
$ws.Range("I11").Value2 = 21410932
$ws.Range("K11").Value2 = 21410932
$ws.Range("M11").Value2 = -21410793
$ws.Range("H33").Value2 = 14597.6
$ws.Range("I33").Value2 = 12552
$ws.Range("K33").Value2 = 12552
$ws.Range("M33").Value2 = -12300
$ws.Range("H63").Value2 = 75498
$ws.Range("I63").Value2 = 39990
$ws.Range("J63").Value2 = 84375
$ws.Range("K63").Value2 = 39990
$ws.Range("L63").Value2 = 84375
$ws.Range("M63").Value2 = -39304
$ws.Range("N63").Value2 = -85747
$ws.Range("H66").Value2 = 75498
$ws.Range("I66").Value2 = 39990
$ws.Range("J66").Value2 = 84375
$ws.Range("K66").Value2 = 119970
$ws.Range("L66").Value2 = 253125
$ws.Range("M66").Value2 = -116538
$ws.Range("N66").Value2 = -259989
$ws.Range("H70").Value2 = 2564762.8
$ws.Range("I70").Value2 = 4837384.5
$ws.Range("J70").Value2 = 8063.375
$ws.Range("K70").Value2 = 4837384.5
$ws.Range("L70").Value2 = 8063.375
$ws.Range("M70").Value2 = -4837114.5
$ws.Range("N70").Value2 = -8603.375
$ws.Range("H73").Value2 = 2564762.8
$ws.Range("I73").Value2 = 4837384.5
$ws.Range("J73").Value2 = 8063.375
$ws.Range("K73").Value2 = 4837384.5
$ws.Range("L73").Value2 = 8063.375
$ws.Range("M73").Value2 = -4836448.5
$ws.Range("N73").Value2 = -9935.375
$ws.Range("H80").Value2 = 1046893.4
$ws.Range("I80").Value2 = 1740755.1
$ws.Range("J80").Value2 = 6100.8335
$ws.Range("K80").Value2 = 1740755.1
$ws.Range("L80").Value2 = 6100.8335
$ws.Range("M80").Value2 = -1739757.1
$ws.Range("N80").Value2 = -8096.8335
$ws.Range("H83").Value2 = 1046893.4
$ws.Range("I83").Value2 = 1740755.1
$ws.Range("J83").Value2 = 6100.8335
$ws.Range("K83").Value2 = 8703775.5
$ws.Range("L83").Value2 = 30504.1675
$ws.Range("M83").Value2 = -8698783.5
$ws.Range("N83").Value2 = -40488.1675
$ws.Range("H102").Value2 = 13843.786
$ws.Range("I102").Value2 = 15402.714
$ws.Range("K102").Value2 = 15402.714
$ws.Range("M102").Value2 = -13780.714
$ws.Range("H122").Value2 = 530168.4
$ws.Range("I122").Value2 = 738902.4
$ws.Range("K122").Value2 = 2216707.2
$ws.Range("M122").Value2 = -2214257.2
$ws.Range("H126").Value2 = 3994.125
$ws.Range("I126").Value2 = 2274
$ws.Range("J126").Value2 = 6027
$ws.Range("K126").Value2 = 6822
$ws.Range("L126").Value2 = 18081
$ws.Range("M126").Value2 = -4352
$ws.Range("N126").Value2 = -23021
$ws.Range("H128").Value2 = 76586.75
$ws.Range("J128").Value2 = 76586.75
$ws.Range("L128").Value2 = 76586.75
$ws.Range("N128").Value2 = -86546.75
$ws.Range("H132").Value2 = 2858.2954
$ws.Range("I132").Value2 = 2826.4324
$ws.Range("J132").Value2 = 3026.7144
$ws.Range("K132").Value2 = 8479.297200000001
$ws.Range("L132").Value2 = 9080.143199999999
$ws.Range("M132").Value2 = -5949.297200000001
$ws.Range("N132").Value2 = -14140.1432
$ws.Range("H141").Value2 = 69999.336
$ws.Range("J141").Value2 = 69999.336
$ws.Range("L141").Value2 = 69999.336
$ws.Range("N141").Value2 = -80359.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 3671.3691
$ws.Range("I7").Value2 = 3082.8394
$ws.Range("K7").Value2 = 3082.8394
$ws.Range("M7").Value2 = -2970.8394
$ws.Range("H16").Value2 = 9029.477000000001
$ws.Range("I16").Value2 = 7632.5625
$ws.Range("K16").Value2 = 7632.5625
$ws.Range("M16").Value2 = -7462.5625
$ws.Range("H22").Value2 = 1584.3529
$ws.Range("I22").Value2 = 1312.3
$ws.Range("J22").Value2 = 1973
$ws.Range("K22").Value2 = 1312.3
$ws.Range("L22").Value2 = 1973
$ws.Range("M22").Value2 = -1017.3
$ws.Range("N22").Value2 = -2563
$ws.Range("H27").Value2 = 1584.3529
$ws.Range("I27").Value2 = 1312.3
$ws.Range("J27").Value2 = 1973
$ws.Range("K27").Value2 = 1312.3
$ws.Range("L27").Value2 = 1973
$ws.Range("M27").Value2 = -1205.3
$ws.Range("N27").Value2 = -2187
$ws.Range("H46").Value2 = 5051.36
$ws.Range("I46").Value2 = 4183.3335
$ws.Range("J46").Value2 = 5169.727
$ws.Range("K46").Value2 = 4183.3335
$ws.Range("L46").Value2 = 5169.727
$ws.Range("M46").Value2 = -3995.3335
$ws.Range("N46").Value2 = -5545.727
$ws.Range("H61").Value2 = 3695.3142
$ws.Range("I61").Value2 = 3639.3447
$ws.Range("J61").Value2 = 3965.8333
$ws.Range("K61").Value2 = 3639.3447
$ws.Range("L61").Value2 = 3965.8333
$ws.Range("M61").Value2 = -3437.3447
$ws.Range("N61").Value2 = -4369.8333
$ws.Range("H93").Value2 = 2516.125
$ws.Range("I93").Value2 = 2447
$ws.Range("K93").Value2 = 2447
$ws.Range("M93").Value2 = -1199
$ws.Range("H113").Value2 = 3695.3142
$ws.Range("I113").Value2 = 3639.3447
$ws.Range("J113").Value2 = 3965.8333
$ws.Range("K113").Value2 = 3639.3447
$ws.Range("L113").Value2 = 3965.8333
$ws.Range("M113").Value2 = -1469.3447
$ws.Range("N113").Value2 = -8305.8333
$ws.Range("H122").Value2 = 7180.909
$ws.Range("I122").Value2 = 3495
$ws.Range("J122").Value2 = 8000
$ws.Range("K122").Value2 = 10485
$ws.Range("L122").Value2 = 24000
$ws.Range("M122").Value2 = -8035
$ws.Range("N122").Value2 = -28900
$ws.Range("H126").Value2 = 3671.3691
$ws.Range("I126").Value2 = 3082.8394
$ws.Range("K126").Value2 = 9248.518199999999
$ws.Range("M126").Value2 = -6778.518199999999
$ws.Range("H129").Value2 = 95000
$ws.Range("J129").Value2 = 95000
$ws.Range("L129").Value2 = 95000
$ws.Range("N129").Value2 = -105000
$ws.Range("H133").Value2 = 100021.75
$ws.Range("J133").Value2 = 100021.75
$ws.Range("L133").Value2 = 100021.75
$ws.Range("N133").Value2 = -105081.75
$ws.Range("H141").Value2 = 45475
$ws.Range("J141").Value2 = 45475
$ws.Range("L141").Value2 = 45475
$ws.Range("N141").Value2 = -55835

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value2 = 16256750
$ws.Range("I13").Value2 = 21671666
$ws.Range("K13").Value2 = 21671666
$ws.Range("M13").Value2 = -21671526
$ws.Range("H132").Value2 = 23818240
$ws.Range("I132").Value2 = 9826.529
$ws.Range("K132").Value2 = 29479.587
$ws.Range("M132").Value2 = -26949.587
$ws.Range("H136").Value2 = 8881.323
$ws.Range("I136").Value2 = 1100
$ws.Range("K136").Value2 = 3300
$ws.Range("M136").Value2 = -750
